# properties to properties tabs of tourney sheets
#
# Adds a new "Properties" sheet (key/value/notes table) to the workbook,
# moving the competition/host/venue key-value pairs that used to live at
# the top/bottom of the "Tournament" sheet into it, along with the four
# palette colors from "Colors".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Trim the "Tournament" sheet: drop the old "competition-key" /
#    "host-key" rows (rows 2:3) and the trailing "venue-key.N" rows
#    (which, after the first delete, land at rows 19:32).
# ---------------------------------------------------------------------
$tournament = $wb.Worksheets.Item("Tournament")
$tournament.Activate()

$tournament.Rows("2:3").Delete()
$tournament.Rows("19:32").Delete()

# restore a sane single-cell selection on the trimmed sheet
$tournament.Range("A15").Select()

# ---------------------------------------------------------------------
# 2. Add the new "Properties" sheet at the end of the workbook.
# ---------------------------------------------------------------------
$sheetCount = $wb.Worksheets.Count
$lastSheet = $wb.Worksheets.Item($sheetCount)
$props = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$props.Name = "Properties"

$rows = @(
  @("key",         "value",                 "notes"),
  @("competition", "mens-copa-america",      ""),
  @("host",        "usa",                    ""),
  @("timezone",    "US/Eastern",             ""),
  @("color.a",     "#b0d0ee",                "pale blue"),
  @("color.b",     "#fab077",                "pale orange"),
  @("color.c",     "#c4e1b5",                "pale green"),
  @("color.d",     "#fee289",                "pale yellow"),
  @("venue.01",    "us-atlanta-ga",          ""),
  @("venue.02",    "us-arlington-tx",        ""),
  @("venue.03",    "us-santa-clara-ca",      ""),
  @("venue.04",    "us-houston-tx",          ""),
  @("venue.05",    "us-miami-fl",            ""),
  @("venue.06",    "us-inglewood-ca",        ""),
  @("venue.07",    "us-kansas-city-ks",      ""),
  @("venue.08",    "us-east-rutherford-nj",  ""),
  @("venue.09",    "us-las-vegas-nv",        ""),
  @("venue.10",    "us-glendale-az",         ""),
  @("venue.11",    "us-orlando-fl",          ""),
  @("venue.12",    "us-austin-tx",           ""),
  @("venue.13",    "us-kansas-city-mo",      ""),
  @("venue.14",    "us-charlotte-nc",        "")
)

for ($r = 0; $r -lt $rows.Length; $r++) {
    $row = $rows[$r]
    for ($c = 0; $c -lt $row.Length; $c++) {
        $value = $row[$c]
        if ($value -ne "") {
            $props.Cells.Item($r + 1, $c + 1).Value = $value
        }
    }
}

# turn the populated range into a table, matching the other "#"-style sheets
$propsRange = $props.Range("A1:C22")
$propsTable = $props.ListObjects.Add(1, $propsRange, [System.Reflection.Missing]::Value, 1)

# approximate the bestFit column widths (ColumnWidth is ~5/6 narrower than
# the raw OOXML column width in this font)
$props.Columns.Item(1).ColumnWidth = 9.83072916666667
$props.Columns.Item(2).ColumnWidth = 16.1666666666667
$props.Columns.Item(3).ColumnWidth = 8.99869791666667

$props.Range("B4").Select()

# make the new sheet the active tab, as in the authored workbook
$props.Activate()
